# Updated cryptos list -- apply the refreshed Price (column D) and
# Volume(1h) (column E) figures scraped for this run onto Sheet1.
#
# The source feed stores every cell as plain text, including Price
# values that happen to look like a number (e.g. "71.73"). Writing such
# a string straight into a General-formatted cell would make Excel
# auto-convert it to a numeric value, so for those specific cells we
# briefly switch the cell to Text format, assign the literal string,
# then restore the cell's style to "Normal" so no stray formatting is
# left behind (matching the original, unformatted cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($address, $text) {
    $cell = $ws.Range($address)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "42.914.69"
$ws.Range("E2").Value = "  +4.69%  "
$ws.Range("D3").Value = "2.275.03"
$ws.Range("E3").Value = "  +5.13%  "
$ws.Range("E4").Value = "  +0.09%  "
Set-TextValue "D5" "250.26"
$ws.Range("E5").Value = "  +1.59%  "
Set-TextValue "D6" "0.636"
$ws.Range("E6").Value = "  +3.99%  "
Set-TextValue "D7" "71.73"
$ws.Range("E7").Value = "  +9.68%  "
$ws.Range("E8").Value = "  -0.14%  "
Set-TextValue "D9" "0.666"
$ws.Range("E9").Value = "  +19.42%  "
Set-TextValue "D10" "39.20"
$ws.Range("E10").Value = "  +12.27%  "
Set-TextValue "D11" "59.71"
$ws.Range("E11").Value = "  +0.36%  "
Set-TextValue "D12" "0.0972"
$ws.Range("E12").Value = "  +5.81%  "
$ws.Range("E13").Value = "  +10.84%  "
Set-TextValue "D14" "0.104"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").Value = "2.614.27"
$ws.Range("E15").Value = "  +5.08%  "
Set-TextValue "D16" "14.91"
$ws.Range("E16").Value = "  +5.37%  "
$ws.Range("E17").Value = "  +5.34%  "
$ws.Range("D18").Value = "2.306.25"
$ws.Range("E18").Value = "  +6.03%  "
$ws.Range("D19").Value = "42.873.99"
$ws.Range("E19").Value = "  +4.89%  "
$ws.Range("E20").Value = "  +8.07%  "
$ws.Range("E21").Value = "  +4.96%  "
Set-TextValue "D22" "73.26"
$ws.Range("E22").Value = "  +3.02%  "
Set-TextValue "D23" "236.20"
$ws.Range("E23").Value = "  +3.67%  "
$ws.Range("E24").Value = "  +3.61%  "
Set-TextValue "D25" "3.97"
$ws.Range("E25").Value = "  +7.28%  "
Set-TextValue "D26" "11.45"
$ws.Range("E26").Value = "  +3.12%  "
Set-TextValue "D27" "0.999"
Set-TextValue "D28" "2.43"
$ws.Range("E28").Value = "  +1.20%  "
$ws.Range("E29").Value = "  -1.29%  "
$ws.Range("E30").Value = "  +9.65%  "
Set-TextValue "D31" "167.75"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("E32").Value = "  +4.97%  "
Set-TextValue "D33" "6.61"
$ws.Range("E33").Value = "  +17.91%  "
$ws.Range("E34").Value = "  +5.80%  "
$ws.Range("E35").Value = "  +9.72%  "
Set-TextValue "D36" "31.18"
$ws.Range("E36").Value = "  +28.87%  "
$ws.Range("E37").Value = "  +4.78%  "
$ws.Range("E38").Value = "  +14.04%  "
$ws.Range("E39").Value = "  +5.58%  "
Set-TextValue "D40" "0.0318"
$ws.Range("E40").Value = "  +6.70%  "
$ws.Range("E41").Value = "  +8.29%  "
$ws.Range("E42").Value = "  +14.26%  "
Set-TextValue "D43" "5.84"
$ws.Range("E43").Value = "  +7.79%  "
$ws.Range("E44").Value = "  +10.33%  "
Set-TextValue "D45" "0.206"
$ws.Range("E45").Value = "  +8.81%  "
Set-TextValue "D46" "62.22"
$ws.Range("E46").Value = "  +4.02%  "
$ws.Range("E47").Value = "  +0.88%  "
$ws.Range("E48").Value = "  +4.19%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  +3.82%  "
$ws.Range("E51").Value = "  +5.20%  "
